$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# New "Wins" / "Race Starts" data for the 20 drivers (rows 2..21),
# in row order as already laid out in column A.
# -----------------------------------------------------------------
$wins   = @(1, 0, 1, 0, 0, 54, 0, 5, 8, 0, 2, 1, 0, 0, 103, 10, 6, 32, 0, 0)
$starts = @(140, 63, 97, 190, 119, 182, 22, 107, 237, 22, 172, 107, 162, 45, 342, 223, 252, 384, 97, 66)

# -----------------------------------------------------------------
# Header row: add the two new column headers
# -----------------------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Height"
$ws.Range("C1").Value = "Weight"
$ws.Range("D1").Value = "Age"
$ws.Range("E1").Value = "Wins"
$ws.Range("F1").Value = "Race Starts"

# -----------------------------------------------------------------
# Write the new Wins / Race Starts columns for every driver row
# -----------------------------------------------------------------
for ($i = 0; $i -lt $wins.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $wins[$i]
    $ws.Cells.Item($row, 6).Value = $starts[$i]
}

# -----------------------------------------------------------------
# Unify formatting: header row and the whole data block now share a
# single plain (non-bold) Arial style, centered + wrapped.
# -----------------------------------------------------------------
$fullRange = $ws.Range("A1:F21")
$fullRange.Borders.LineStyle = -4142
$fullRange.Font.Name = "Arial"
$fullRange.Font.Size = 11
$fullRange.Font.Bold = $false
$fullRange.Font.Color = 1907739
$fullRange.HorizontalAlignment = -4108
$fullRange.VerticalAlignment = -4108
$fullRange.WrapText = $true

# Header row is taller to match the new, bigger column set
$ws.Rows.Item(1).RowHeight = 28

# -----------------------------------------------------------------
# View state: scroll back to the top and select A2
# -----------------------------------------------------------------
[void]$ws.Range("A2").Select()

Write-Output "done"
